# Week 19 profile update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark this week's submitters (column C) as done: 0 -> 1
$ws.Range("C28").Value = 1
$ws.Range("C54").Value = 1
$ws.Range("C63").Value = 1
$ws.Range("C85").Value = 1
$ws.Range("C90").Value = 1

# Restore the editor's scroll position / active selection from when the
# workbook was saved (row 35 at the top, cell H56 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$ws.Range("H56").Select()
